$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --------------------------------------------------------------------------
# Cache the existing cell-formats we will need onto scratch cells (row 50+)
# *before* any destructive edit, since several source cells below get
# overwritten with new values later in this script.
#   old style 9  (plain border-left cell)                         -> A50
#   old style 15 (fontId3 borderId5, vertical-top alignment)      -> A51
#   old style 12 (fontId3 borderId0)                              -> A52
#   old style 11 (fontId3 borderId5)                              -> A53
#   old style 13 (fontId3 numFmt2 borderId0)                      -> A54
#   old style 14 (fontId4 numFmt2 borderId0)                      -> A55
#   old style 16 (fontId3 borderId2, horizontal-left, inlineStr)  -> A56
#   old style 17 (fontId3 fillId0 borderId2, horizontal-left)     -> A57
# --------------------------------------------------------------------------
$ws.Range("H5").Copy()
$ws.Range("A50").PasteSpecial(-4122)
$ws.Range("A6").Copy()
$ws.Range("A51").PasteSpecial(-4122)
$ws.Range("B5").Copy()
$ws.Range("A52").PasteSpecial(-4122)
$ws.Range("A5").Copy()
$ws.Range("A53").PasteSpecial(-4122)
$ws.Range("D5").Copy()
$ws.Range("A54").PasteSpecial(-4122)
$ws.Range("F5").Copy()
$ws.Range("A55").PasteSpecial(-4122)
$ws.Range("A9").Copy()
$ws.Range("A56").PasteSpecial(-4122)
$ws.Range("B9").Copy()
$ws.Range("A57").PasteSpecial(-4122)

# --------------------------------------------------------------------------
# Insert a new row at 9 -- this pushes the old footer row (9) down to row
# 10 and updates the A9:G9 merge to A10:G10 automatically.
# --------------------------------------------------------------------------
$ws.Rows("9").Insert()

# --------------------------------------------------------------------------
# Re-point the "Cylinder" merge: was A6:A7, now split into A5:A6 / A7:A8.
# --------------------------------------------------------------------------
$ws.Range("A6:A7").UnMerge()
$ws.Range("A5:A6").Merge()
$ws.Range("A7:A8").Merge()

# --------------------------------------------------------------------------
# Row 5 (4-cylinder, automatic transmission)
# --------------------------------------------------------------------------
$ws.Range("A51").Copy(); $ws.Range("A5").PasteSpecial(-4122)
$ws.Range("A5").Value = 4
$ws.Range("A52").Copy(); $ws.Range("B5").PasteSpecial(-4122)
$ws.Range("B5").Value = 0
$ws.Range("A53").Copy(); $ws.Range("C5").PasteSpecial(-4122)
$ws.Range("C5").Value = 1
$ws.Range("A54").Copy(); $ws.Range("D5").PasteSpecial(-4122)
$ws.Range("D5").Value = 91
$ws.Range("A54").Copy(); $ws.Range("E5").PasteSpecial(-4122)
$ws.Range("E5").ClearContents()
$ws.Range("A55").Copy(); $ws.Range("F5").PasteSpecial(-4122)
$ws.Range("F5").Value = 2.14
$ws.Range("A55").Copy(); $ws.Range("G5").PasteSpecial(-4122)
$ws.Range("G5").ClearContents()

# --------------------------------------------------------------------------
# Row 6 (4-cylinder, manual transmission)
# --------------------------------------------------------------------------
$ws.Range("A50").Copy(); $ws.Range("A6").PasteSpecial(-4122)
$ws.Range("A6").ClearContents()
$ws.Range("A52").Copy(); $ws.Range("B6").PasteSpecial(-4122)
$ws.Range("B6").Value = 1
$ws.Range("A53").Copy(); $ws.Range("C6").PasteSpecial(-4122)
$ws.Range("C6").Value = 10
$ws.Range("A55").Copy(); $ws.Range("D6").PasteSpecial(-4122)
$ws.Range("D6").Value = 81.8
$ws.Range("A54").Copy(); $ws.Range("E6").PasteSpecial(-4122)
$ws.Range("E6").Value = 21.87235698318771
$ws.Range("A54").Copy(); $ws.Range("F6").PasteSpecial(-4122)
$ws.Range("F6").Value = 2.3003
$ws.Range("A54").Copy(); $ws.Range("G6").PasteSpecial(-4122)
$ws.Range("G6").Value = 0.5982073312080948

# --------------------------------------------------------------------------
# Row 7 (6-cylinder, automatic transmission)
# --------------------------------------------------------------------------
$ws.Range("A51").Copy(); $ws.Range("A7").PasteSpecial(-4122)
$ws.Range("A7").Value = 6
$ws.Range("A52").Copy(); $ws.Range("B7").PasteSpecial(-4122)
$ws.Range("B7").Value = 0
$ws.Range("A53").Copy(); $ws.Range("C7").PasteSpecial(-4122)
$ws.Range("C7").Value = 3
$ws.Range("A55").Copy(); $ws.Range("D7").PasteSpecial(-4122)
$ws.Range("D7").Value = 131.6666666666667
$ws.Range("A54").Copy(); $ws.Range("E7").PasteSpecial(-4122)
$ws.Range("E7").Value = 37.52776749732568
$ws.Range("A54").Copy(); $ws.Range("F7").PasteSpecial(-4122)
$ws.Range("F7").Value = 2.755
$ws.Range("A54").Copy(); $ws.Range("G7").PasteSpecial(-4122)
$ws.Range("G7").Value = 0.1281600561797629

# --------------------------------------------------------------------------
# Row 8 (6-cylinder, manual transmission)
# --------------------------------------------------------------------------
$ws.Range("A50").Copy(); $ws.Range("A8").PasteSpecial(-4122)
$ws.Range("A8").ClearContents()
$ws.Range("A52").Copy(); $ws.Range("B8").PasteSpecial(-4122)
$ws.Range("B8").Value = 1
$ws.Range("A53").Copy(); $ws.Range("C8").PasteSpecial(-4122)
$ws.Range("C8").Value = 4
$ws.Range("A54").Copy(); $ws.Range("D8").PasteSpecial(-4122)
$ws.Range("D8").Value = 115.25
$ws.Range("A54").Copy(); $ws.Range("E8").PasteSpecial(-4122)
$ws.Range("E8").Value = 9.178779875342908
$ws.Range("A55").Copy(); $ws.Range("F8").PasteSpecial(-4122)
$ws.Range("F8").Value = 3.38875
$ws.Range("A55").Copy(); $ws.Range("G8").PasteSpecial(-4122)
$ws.Range("G8").Value = 0.1162163929916946

# --------------------------------------------------------------------------
# Row 9 (8-cylinder, all transmissions combined) -- brand-new row
# --------------------------------------------------------------------------
$ws.Range("A53").Copy(); $ws.Range("A9").PasteSpecial(-4122)
$ws.Range("A9").Value = 8
$ws.Range("A52").Copy(); $ws.Range("B9").PasteSpecial(-4122)
$ws.Range("B9").Value = 0
$ws.Range("A53").Copy(); $ws.Range("C9").PasteSpecial(-4122)
$ws.Range("C9").Value = 14
$ws.Range("A54").Copy(); $ws.Range("D9").PasteSpecial(-4122)
$ws.Range("D9").Value = 209.2142857142857
$ws.Range("A54").Copy(); $ws.Range("E9").PasteSpecial(-4122)
$ws.Range("E9").Value = 50.97688551827051
$ws.Range("A54").Copy(); $ws.Range("F9").PasteSpecial(-4122)
$ws.Range("F9").Value = 3.999214285714287
$ws.Range("A54").Copy(); $ws.Range("G9").PasteSpecial(-4122)
$ws.Range("G9").Value = 0.7594047444769265
$ws.Range("A50").Copy(); $ws.Range("H9").PasteSpecial(-4122)
$ws.Range("H9").ClearContents()

# --------------------------------------------------------------------------
# Clean up the scratch cells used for format-caching.
# --------------------------------------------------------------------------
$ws.Range("A50:A57").Clear()

$ws.Range("A1").Select()
